$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")

# Total IN
$summary.Range("B3").Value = 2

# Gender Breakdown: Female -> Unknown
$summary.Range("A10").Value = "Unknown"
$summary.Range("B10").Value = 1

# Age Group Breakdown: 13-25 row
$summary.Range("B16").Value = 1
$summary.Range("C16").Value = 1

# Age Group Breakdown: 26-45 row
$summary.Range("B17").Value = 0
$summary.Range("C17").Value = 0

# --- Sheet: Hourly Breakdown ---
$hourly = $wb.Worksheets.Item("Hourly Breakdown")

# 14:00 row
$hourly.Range("B16").Value = 2
$hourly.Range("D16").Value = 0

# --- Sheet: Charts Data ---
$charts = $wb.Worksheets.Item("Charts Data")

# Age Distribution (IN): 26-45 -> 13-25
$charts.Range("A8").Value = "13-25"
